$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank "=" separator row above the table header (old row 9 -> new row 10, etc.)
$ws.Rows.Item(9).Insert()

# Populate the new row 9 by copying the still-unedited "=" marker row (row 8) down,
# avoiding the leading "=" being interpreted as a formula if assigned via .Value directly.
$ws.Range("A8:G8").Copy($ws.Range("A9:G9"))

# Update the three "test" rows: Ignore test is new, pushing Merge/Append down one slot.
$ws.Range("B6").Value = "Ignore test:"
$ws.Range("C6").Value = "xltablediff.py  --key ID --ignore Color test1old.xlsx test1new.xlsx --out test1ignore.xlsx"

$ws.Range("B7").Value = "Merge test:"
$ws.Range("C7").Value = "xltablediff.py  --key ID --merge Color test1old.xlsx test1new.xlsx --out test1merge.xlsx"

$ws.Range("B8").Value = "Append test:"
$ws.Range("C8").Value = "xltablediff.py  --key ID --append test1old.xlsx test1new.xlsx --out test1append.xlsx"
